# "Scripts updated as per failures" - the automation test data (account
# sign-up credentials the Selenium/CreateAccount sheet feeds into the web
# form) was refreshed after a test re-run. Each of these cells had already
# been stamped with a "filled" look (white fill + thin top/bottom rule) by
# the framework; we redo the value + restamp the same look so the sheet
# ends up in the same state the authoring tool produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateAccount")

function Stamp-Field($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.Value = $value
    $rng.Interior.ColorIndex = 2
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(9).LineStyle = 1
}

Stamp-Field "E2" "SeleniumGxjc@mailinator.com"
Stamp-Field "E3" "SeleniumfJAa@mailinator.com"
Stamp-Field "E4" "SeleniumUORD@mailinator.com"
Stamp-Field "E5" "SeleniumGjBp@mailinator.com"
Stamp-Field "E6" "SeleniumqoDh@mailinator.com"
Stamp-Field "F7" "Automation8509!"

# Column E keeps "best fit" (it was already best-fit) now that the
# refreshed e-mail addresses are in place.
$ws.Columns.Item(5).AutoFit()
